$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date in C1 (3/15/2024 -> 3/28/2024) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "RAF-capacity" sheet: raise the capacity-credit multiplier for the two
#     hydrogen rows (B24 "hydrogen combustion turbine", B25 "hydrogen combined
#     cycle") from 0.3 up to 1, widen column A so the longer labels are
#     readable, and update the selection/scroll position ---
$wsCap = $wb.Worksheets.Item("RAF-capacity")
$wsCap.Range("B24").Value = 1
$wsCap.Range("B25").Value = 1
$wsCap.Columns.Item(1).ColumnWidth = 28.2
$wsCap.Range("B25").Select()

# --- Switch the active/selected tab from "RAF-generation" to "RAF-capacity",
#     and bump the zoom level on that sheet to 80% ---
$wsCap.Activate()
$excel.ActiveWindow.Zoom = 80
